$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "averageDIAG": convert A and C columns to shared formulas,
# add a new E1 average cell, and drop the old row 46 (average moved
# up into row 1).
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("averageDIAG")

$ws1.Range("A1:A45").Formula = "=1/COS(C1)"
$ws1.Range("C1:C45").Formula = "=RADIANS(B1)"
$ws1.Range("E1").Formula = "=AVERAGE(A1:A45)"
$ws1.Rows.Item(46).Delete()

# ------------------------------------------------------------------
# Sheet "Computing": rebuilt layout (E0 correction block in rows 1-3,
# PD block in rows 7-10), with "Good"/"Neutral" highlight styles and
# a bold header style.
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Computing")

# Wipe the old layout first so stale shared strings / formulas don't linger.
$ws2.Range("A1:K10").Clear()

# --- text labels, entered in the same order the workbook's shared
#     string table lists them so the rebuilt table lines up exactly.
$ws2.Range("A1").Value = "E0"
$ws2.Range("J2").Value = "PD"
$ws2.Range("K2").Value = "pond"
$ws2.Range("J3").Value = "ShapeFactor"
$ws2.Range("F1").Value = "E0_mm"
$ws2.Range("C1").Value = "PD_w_m"
$ws2.Range("D1").Value = "PDarea_m2"
$ws2.Range("C7").Value = "PDw_m"
$ws2.Range("D7").Value = "PDh_m"
$ws2.Range("B1").Value = "cell size_m"
$ws2.Range("E1").Value = "PDarea_ha"
$ws2.Range("G1").Value = "E0corr_mm"
$ws2.Range("H1").Value = "E0corr_mm_FORMULA"
$ws2.Range("E7").Value = "PDvol_m3"
$ws2.Range("G7").Value = "PD_mm"
$ws2.Range("F7").Value = "PDvol_m3_ha"
$ws2.Range("H7").Value = "PD_mm_FORMULA"
$ws2.Range("A7").Value = "PD"
$ws2.Range("B7").Value = "cell size_m"

# --- ShapeFactor helper (averages the 1/cos(lat) table on the other sheet)
$ws2.Range("K3").Formula = "=averageDIAG!E1"

# --- E0 correction block --------------------------------------------------
$ws2.Range("B2").Value = 100
$ws2.Range("C2").Value = 1
$ws2.Range("F2").Value = 3
$ws2.Range("D2").Formula = "=C2*`$K`$3*B2"
$ws2.Range("E2").Formula = "=D2/100^2"
$ws2.Range("G2").Formula = "=E2*F2"
$ws2.Range("H2").Formula = "=B2*C2*`$K`$3*F2/100^2"

$ws2.Range("B3").Value = 10
$ws2.Range("C3").Value = 1
$ws2.Range("F3").Value = 3
$ws2.Range("D3").Formula = "=C3*`$K`$3*B3"
$ws2.Range("E3").Formula = "=D3/100^2"
$ws2.Range("G3").Formula = "=E3*F3"
$ws2.Range("H3").Formula = "=B3*C3*`$K`$3*F3/100^2"

# --- PD (pond) block -------------------------------------------------------
$ws2.Range("B8").Value = 100
$ws2.Range("C8").Value = 1
$ws2.Range("D8").Value = 2
$ws2.Range("E8").Formula = "=C8*D8*B8*`$K`$3"
$ws2.Range("F8").Formula = "=E8*B8^2/100^2"
$ws2.Range("G8").Formula = "=F8/10"
$ws2.Range("H8").Formula = "=C8*D8*B8*`$K`$3*B8^2/100^2/10"

$ws2.Range("B9").Value = 10
$ws2.Range("C9").Value = 1
$ws2.Range("D9").Value = 2
$ws2.Range("B10").Value = 50
$ws2.Range("C10").Value = 1
$ws2.Range("D10").Value = 2

$ws2.Range("E9:E10").Formula = "=C9*D9*B9*`$K`$3"
$ws2.Range("F9:F10").Formula = "=E9*B9^2/100^2"
$ws2.Range("G9:G10").Formula = "=F9/10"
$ws2.Range("H9:H10").Formula = "=C9*D9*B9*`$K`$3*B9^2/100^2/10"

# --- styling ----------------------------------------------------------------
$ws2.Range("A1").Font.Bold = $true
$ws2.Range("A7").Font.Bold = $true

$ws2.Range("G1:H3").Style = "Good"
$ws2.Range("G7:H10").Style = "Good"

$ws2.Range("B2:C3").Style = "Neutral"
$ws2.Range("B8:D10").Style = "Neutral"

$ns = $wb.Styles.Item("Neutral")
$ns.Font.Color = 26012

# --- column widths & selection ----------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 3.42578125
$ws2.Columns.Item(2).ColumnWidth = 10.7109375
$ws2.Columns.Item(3).ColumnWidth = 8.7109375
$ws2.Columns.Item(4).ColumnWidth = 12
$ws2.Columns.Item(5).ColumnWidth = 12
$ws2.Columns.Item(6).ColumnWidth = 13.140625
$ws2.Columns.Item(7).ColumnWidth = 12
$ws2.Columns.Item(8).ColumnWidth = 20.85546875
$ws2.Columns.Item(10).ColumnWidth = 11.85546875
$ws2.Columns.Item(11).ColumnWidth = 12

$ws2.Range("H9").Select()
